# Update "想去人数" (F column) counts across the four sheets to match the
# refreshed data snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1427
$ws.Range("F3").Value = 105
$ws.Range("F4").Value = 2109
$ws.Range("F5").Value = 6774
$ws.Range("F6").Value = 530
$ws.Range("F7").Value = 1059
$ws.Range("F9").Value = 4585
$ws.Range("F10").Value = 6807
$ws.Range("F12").Value = 228
$ws.Range("F14").Value = 812
$ws.Range("F15").Value = 119
$ws.Range("F18").Value = 1132
$ws.Range("F20").Value = 132
$ws.Range("F24").Value = 1065
$ws.Range("F28").Value = 121
$ws.Range("F30").Value = 1165
$ws.Range("F34").Value = 10
$ws.Range("F37").Value = 516
$ws.Range("F38").Value = 370
$ws.Range("F39").Value = 43
$ws.Range("F41").Value = 316
$ws.Range("F43").Value = 526

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 517
$ws.Range("F22").Value = 193
$ws.Range("F31").Value = 792
$ws.Range("F32").Value = 967

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 612
$ws.Range("F8").Value = 1304
$ws.Range("F9").Value = 1950

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1427
$ws.Range("F7").Value = 105
$ws.Range("F8").Value = 612
$ws.Range("F9").Value = 612
$ws.Range("F10").Value = 517
$ws.Range("F11").Value = 6774
$ws.Range("F12").Value = 530
$ws.Range("F13").Value = 1059
$ws.Range("F15").Value = 4585
$ws.Range("F17").Value = 6807
$ws.Range("F19").Value = 228
$ws.Range("F22").Value = 812
$ws.Range("F23").Value = 119
$ws.Range("F24").Value = 1304
$ws.Range("F25").Value = 1950
$ws.Range("F26").Value = 193
$ws.Range("F27").Value = 1132
$ws.Range("F28").Value = 132
$ws.Range("F30").Value = 1065
$ws.Range("F33").Value = 121
$ws.Range("F34").Value = 1165
$ws.Range("F38").Value = 10
$ws.Range("F39").Value = 967
$ws.Range("F40").Value = 516
$ws.Range("F42").Value = 370
$ws.Range("F43").Value = 43
$ws.Range("F45").Value = 316
$ws.Range("F46").Value = 526
